$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the accessory insured-sum label for "Movilidad"
$ws.Range("B2").Value = "Hasta $150.000"

# Move the active selection to B3
$ws.Range("B3").Select()
